# Refresh the screener table on Sheet1 with the new ticker lists.
# Columns: A=index, B=Buying Opportunity, C=support Zone,
#          D=long buildup, E=Short buildup, F=FII ENTERING.
# The table grows from 15 data rows (2-16) to 28 data rows (2-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Grow the table: rows 2-16 already carry the bordered/bold/centered
#    style (s="1") on column A. Copy that formatting down to the 13 new
#    rows (17-29) before writing their values, so the new index cells
#    pick up the same style instead of minting a new one.
$ws.Range("A16").Copy($ws.Range("A17:A29"))

# 2) Write the full target content for every data row.
$rows = @(
  @{ Row=2;  A=0;  B="NSE:ALOKINDS";   C="NSE:ADANIGREEN"; D="NSE:ATUL"; E="NSE:MANAPPURAM"; F="NSE:BALKRISIND" },
  @{ Row=3;  A=1;  B="NSE:ANGELONE";   C="NSE:ADORWELD";   D=""; E=""; F="" },
  @{ Row=4;  A=2;  B="NSE:BALKRISIND"; C="NSE:ASHIANA";    D=""; E=""; F="" },
  @{ Row=5;  A=3;  B="NSE:BODALCHEM";  C="NSE:ASKAUTOLTD"; D=""; E=""; F="" },
  @{ Row=6;  A=4;  B="NSE:DHANUKA";    C="NSE:CAMPUS";     D=""; E=""; F="" },
  @{ Row=7;  A=5;  B="NSE:DHUNINV";    C="NSE:CAPACITE";   D=""; E=""; F="" },
  @{ Row=8;  A=6;  B="NSE:EROSMEDIA";  C="NSE:CREDITACC";  D=""; E=""; F="" },
  @{ Row=9;  A=7;  B="NSE:IIFLSEC";    C="NSE:EASEMYTRIP"; D=""; E=""; F="" },
  @{ Row=10; A=8;  B="NSE:IMFA";       C="NSE:EQUITASBNK"; D=""; E=""; F="" },
  @{ Row=11; A=9;  B="NSE:INDIAMART";  C="NSE:GLOBUSSPR";  D=""; E=""; F="" },
  @{ Row=12; A=10; B="NSE:LFIC";       C="NSE:HARDWYN";    D=""; E=""; F="" },
  @{ Row=13; A=11; B="NSE:MANGCHEFER"; C="NSE:HESTERBIO";  D=""; E=""; F="" },
  @{ Row=14; A=12; B="NSE:MCLEODRUSS"; C="NSE:HINDUNILVR"; D=""; E=""; F="" },
  @{ Row=15; A=13; B="NSE:PRAKASHSTL"; C="NSE:ITC";        D=""; E=""; F="" },
  @{ Row=16; A=14; B="NSE:RPGLIFE";    C="NSE:JSWENERGY";  D=""; E=""; F="" },
  @{ Row=17; A=15; B="";               C="NSE:JUBLINDS";   D=""; E=""; F="" },
  @{ Row=18; A=16; B="";               C="NSE:LODHA";      D=""; E=""; F="" },
  @{ Row=19; A=17; B="";               C="NSE:LT";         D=""; E=""; F="" },
  @{ Row=20; A=18; B="";               C="NSE:MAHLIFE";    D=""; E=""; F="" },
  @{ Row=21; A=19; B="";               C="NSE:NAHARCAP";   D=""; E=""; F="" },
  @{ Row=22; A=20; B="";               C="NSE:NUCLEUS";    D=""; E=""; F="" },
  @{ Row=23; A=21; B="";               C="NSE:PAVNAIND";   D=""; E=""; F="" },
  @{ Row=24; A=22; B="";               C="NSE:POWERGRID";  D=""; E=""; F="" },
  @{ Row=25; A=23; B="";               C="NSE:PRSMJOHNSN"; D=""; E=""; F="" },
  @{ Row=26; A=24; B="";               C="NSE:RBA";        D=""; E=""; F="" },
  @{ Row=27; A=25; B="";               C="NSE:REPL";       D=""; E=""; F="" },
  @{ Row=28; A=26; B="";               C="NSE:RITES";      D=""; E=""; F="" },
  @{ Row=29; A=27; B="";               C="NSE:RTNINDIA";   D=""; E=""; F="" }
)

foreach ($r in $rows) {
  $ws.Range("A" + $r.Row).Value = $r.A
  $ws.Range("B" + $r.Row).Value = $r.B
  $ws.Range("C" + $r.Row).Value = $r.C
  $ws.Range("D" + $r.Row).Value = $r.D
  $ws.Range("E" + $r.Row).Value = $r.E
  $ws.Range("F" + $r.Row).Value = $r.F
}

Write-Host "Updated $($rows.Count) data rows (2-29) on $($ws.Name)."
